# Add the student's name as a new first paragraph, followed by a blank
# paragraph, ahead of the existing document content.
$d = $word.ActiveDocument

# Insert a paragraph mark right before the very start of the document,
# then another one before that — this yields two new empty paragraphs
# in front of the original first paragraph ("Milestones:").
$r = $d.Range(0, 0)
$r.InsertParagraphBefore()

$r2 = $d.Range(0, 0)
$r2.InsertParagraphBefore()

# Fill the first of the two new paragraphs with the student's name;
# the second stays empty, matching the target layout.
$p1 = $d.Paragraphs(1)
$p1.Range.Text = "Aluno: João Vitor da Silva Neto"
